# Remove the redundant "BLM-RB4-" prefix from the StationCode values (column A)
# on every worksheet, and narrow column A's width to fit the shorter text.

$wb = $excel.ActiveWorkbook
$prefix = "BLM-RB4-"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count
    $firstRowCell = $used.Cells.Item(1, 1)
    $startRow = $firstRowCell.Row

    for ($r = $startRow; $r -lt ($startRow + $lastRow); $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $val = $cell.Value2
        if ($val -ne $null -and $val -like "$prefix*") {
            $cell.Value2 = $val.Substring($prefix.Length)
        }
    }

    # Column width was 18 characters wide; shrink it to 13 now that the
    # "BLM-RB4-" prefix has been stripped from every station code.
    $ws.Columns.Item(1).ColumnWidth = 12.14
}
